$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing bold/border/centered style (currently on column-A label
# cells A7:A11) before we touch anything, then stamp it onto the new header
# row (B1:E1) and the new label column (A2:A3) via copy/paste-special so we
# reuse the existing style index instead of minting new ones.
$ws.Range("A9").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Wipe out the old table (rows 7-11) entirely.
$ws.Range("A7:E11").Clear()

# New header row.
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"
$ws.Range("E1").Value = "Unnamed: 4"

# New data rows (column by column, matching the authored write order).
$ws.Range("A2").Value = "2d0f4fde560a87908cb87b7a0a29cebd00c9f7fd38d08dbedb41ce22dd7ad879:bfb76a73f68f4327afdf0b201caf51a2"
$ws.Range("A3").Value = "30cd31772a084fc1a31447c7e845aa4b1828ad34091d74393ff6277a66c316f7:18bc1c39b86c43e2825353c3b05c5f9a"

$ws.Range("B2").Value = "hlkhljkhl"
$ws.Range("B3").Value = "tuituit"

$ws.Range("C2").Value = "hjklhljh"
$ws.Range("C3").Value = "tuiti"

$ws.Range("D2").Value = "jhlkh"
$ws.Range("D3").Value = "tuiutiutiti"

$ws.Range("E2").Value = "hjklkhhhlklh"
$ws.Range("E3").Value = "tuit"

# Match the author's final selection.
$ws.Range("F17").Select()
